$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.032280325889587
$ws.Range("B1").Value = 1.917413949966431
$ws.Range("C1").Value = 7.964878082275391
$ws.Range("D1").Value = 1.999763369560242
$ws.Range("E1").Value = 0.6303097605705261
